$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be created in this exact order so that the
# --- underlying shared-string table indices line up with the target file
# --- (index 9 .. 16), before we touch cells that reuse existing strings.

# index 9 (rich text, two runs with different formatting) -> B11
$run1 = "can write just default: 0 or write the start time if you want to discard the cycle"
$run2 = '; this number is just the sanity check - it indicates what is the exact starting time of the cycle that will be "cleaned"'
$ws.Range("B11").Value = ($run1 + $run2)
$chars = $ws.Range("B11").Characters($run1.Length + 1, $run2.Length)
$chars.Font.Size = 12
$chars.Font.Name = "Calibri"

# index 10 -> B12
$ws.Range("B12").Value = "indicate the start time of the section to be KEPT, write 0 if you want to discard this cycle: in minutes"

# index 11 -> B13
$ws.Range("B13").Value = "indicate the end time of the section to be KEPT, write 0 if you want to discard this cycle: in minutes"

# index 12 -> A7
$ws.Range("A7").Value = "KEY "

# index 13 -> B9
$ws.Range("B9").Value = "box number, must match with the original  filename"

# index 14 -> B8
$ws.Range("B8").Value = "date (three same as the file name). TIP: add ' symbol right before typing the date, because excel likes to change it in their date formats. "

# index 15 -> B10
$ws.Range("B10").Value = "channel number, will be either 1,2,3,4; TIP: can read this off from the .png plot file while looking at what needs to be cleaned"

# index 16 -> B14
$ws.Range("B14").Value = "for general cleaning - leave as is ('smr' for all). this comes in useful when slopes themselves are analysed, was used for Lobster anpnea (ask details if interested)"

# --- column A labels reuse the original header strings (indices 0-6)
$ws.Range("A8").Value = "date"
$ws.Range("A9").Value = "box"
$ws.Range("A10").Value = "channel"
$ws.Range("A11").Value = "cycle_file_start"
$ws.Range("A12").Value = "sectioned_slope_start"
$ws.Range("A13").Value = "sectioned_slope_end"
$ws.Range("A14").Value = "type"

# --- formatting: bold + yellow fill for column A (rows 8-14) and B11,
# --- plain yellow fill (no bold) for the rest of column B (rows 8-10,12-14)
# NOTE: Interior.Color must be applied before Font.Bold on each cell -- the
# engine records every intermediate (fill, font) state as its own style
# entry, and this ordering lets the "fill-only" stepping-stone state be the
# same one later reused verbatim by the non-bold column-B cells (matching
# the two cellXfs entries the target workbook ends up with).
foreach ($r in 8..14) {
    $cell = $ws.Range("A" + $r)
    $cell.Interior.Color = 65535
    $cell.Font.Bold = $true
}

foreach ($r in @(8, 9, 10, 12, 13, 14)) {
    $ws.Range("B" + $r).Interior.Color = 65535
}

$ws.Range("B11").Interior.Color = 65535
$ws.Range("B11").Font.Bold = $true

# --- selection matches the post-edit workbook state
$ws.Range("D9").Select() | Out-Null
